$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Lesson Names" table (B3:C..) is sorted alphabetically by the
# "Original Name" column. Two new lesson-name mappings were added to the
# source data, which land in their correct alphabetical position and push
# every following row down by one.

# 1) New row: "Architecture et urbanisme" / "Architecture"
#    Alphabetically this sits between "Anglais" (row 8) and "Arts" (row 9),
#    so insert a new row at row 9.
$ws.Rows(9).Insert()
$ws.Range("B9").Value = "Architecture et urbanisme"
$ws.Range("C9").Value = "Architecture"

# 2) New row: "Géographie territoriale appliquée" / "Géographie territoriale"
#    Alphabetically this sits between "Géographie" and "Géométrie descriptive".
#    After the first insertion above, "Géométrie descriptive" now lives at
#    row 45, so the new row is inserted at row 44.
$ws.Rows(44).Insert()
$ws.Range("B44").Value = "Géographie territoriale appliquée"
$ws.Range("C44").Value = "Géographie territoriale"

# Update the active selection to match the saved view state.
$ws.Range("D7").Select()
